$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2446.0688
$ws.Range("I40").Value = 1775.3334
$ws.Range("J40").Value = 2621.0435
$ws.Range("K40").Value = 1775.3334
$ws.Range("L40").Value = 2621.0435
$ws.Range("M40").Value = -1600.3334
$ws.Range("N40").Value = -2971.0435

# Row 137
$ws.Range("H137").Value = 1667.8286
$ws.Range("I137").Value = 1684.0605
$ws.Range("J137").Value = 1400
$ws.Range("K137").Value = 5052.181500000001
$ws.Range("L137").Value = 4200
$ws.Range("M137").Value = -2502.181500000001
$ws.Range("N137").Value = -9300


$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2909.9656
$ws.Range("I45").Value = 1344
$ws.Range("J45").Value = 4182.3125
$ws.Range("K45").Value = 1344
$ws.Range("L45").Value = 4182.3125
$ws.Range("M45").Value = -967
$ws.Range("N45").Value = -4936.3125

# Row 98
$ws.Range("H98").Value = 27666.666
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 27666.666
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 27666.666
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -33656.666

# Row 132
$ws.Range("H132").Value = 8910
$ws.Range("I132").Value = 11380.8
$ws.Range("J132").Value = 8223.666999999999
$ws.Range("K132").Value = 34142.39999999999
$ws.Range("L132").Value = 24671.001
$ws.Range("M132").Value = -31612.39999999999
$ws.Range("N132").Value = -29731.001


$ws = $wb.Worksheets.Item("BSM")
# Row 88
$ws.Range("H88").Value = 22000.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 22000.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 22000.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -22812.5

# Row 91
$ws.Range("H91").Value = 22000.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 22000.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 22000.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -24808.5

# Row 102
$ws.Range("H102").Value = 8556
$ws.Range("I102").Value = 8556
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8556
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -5311
$ws.Range("N102").ClearContents()

# Row 134
$ws.Range("H134").Value = 3317.919
$ws.Range("I134").Value = 1446.1041
$ws.Range("J134").Value = 6773.577
$ws.Range("K134").Value = 4338.3123
$ws.Range("L134").Value = 20320.731
$ws.Range("M134").Value = -1803.3123
$ws.Range("N134").Value = -25390.731


$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 2337
$ws.Range("I33").Value = 1015.5
$ws.Range("J33").Value = 4980
$ws.Range("K33").Value = 1015.5
$ws.Range("L33").Value = 4980
$ws.Range("M33").Value = -636.5
$ws.Range("N33").Value = -5738

# Row 41
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -20856

# Row 42
$ws.Range("H42").Value = 8000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 8000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 8000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -9186

# Row 50
$ws.Range("H50").Value = 20000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20000
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -21250

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# Row 55
$ws.Range("H55").Value = 10900
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 10900
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10900
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -11530

# Row 60
$ws.Range("H60").Value = 9093
$ws.Range("I60").Value = 9093
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 9093
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8582
$ws.Range("N60").ClearContents()

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# Row 74
$ws.Range("H74").Value = 33624.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 33624.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 33624.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -35372.5

# Row 77
$ws.Range("H77").Value = 33624.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 33624.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 100873.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -109609.5

# Row 88
$ws.Range("H88").Value = 25448.6
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 25448.6
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 25448.6
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -26260.6

# Row 91
$ws.Range("H91").Value = 25448.6
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 25448.6
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 25448.6
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -28256.6


$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 367.6154
$ws.Range("I17").Value = 336.27274
$ws.Range("J17").Value = 540
$ws.Range("K17").Value = 1008.81822
$ws.Range("L17").Value = 1620
$ws.Range("M17").Value = -839.81822
$ws.Range("N17").Value = -1958

# Row 20
$ws.Range("H20").Value = 2390.6
$ws.Range("I20").Value = 2035.3103
$ws.Range("J20").Value = 3327.2727
$ws.Range("K20").Value = 6105.9309
$ws.Range("L20").Value = 9981.8181
$ws.Range("M20").Value = -5878.9309
$ws.Range("N20").Value = -10435.8181

# Row 21
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -2827
$ws.Range("N21").ClearContents()

# Row 26
$ws.Range("H26").Value = 109.42857
$ws.Range("I26").Value = 79.53333000000001
$ws.Range("J26").Value = 184.16667
$ws.Range("K26").Value = 238.59999
$ws.Range("L26").Value = 552.50001
$ws.Range("M26").Value = 49.40000999999998
$ws.Range("N26").Value = -1128.50001

# Row 131
$ws.Range("H131").Value = 1106.1428
$ws.Range("I131").Value = 826.6667
$ws.Range("J131").Value = 1152.7222
$ws.Range("K131").Value = 2480.0001
$ws.Range("L131").Value = 3458.1666
$ws.Range("M131").Value = 2559.9999
$ws.Range("N131").Value = -13538.1666

# Row 133
$ws.Range("H133").Value = 3146
$ws.Range("I133").Value = 1976.6666
$ws.Range("J133").Value = 4900
$ws.Range("K133").Value = 5929.9998
$ws.Range("L133").Value = 14700
$ws.Range("M133").Value = -869.9997999999996
$ws.Range("N133").Value = -24820


$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 11774.667
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 11774.667
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 11774.667
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -17266.667

# Row 113
$ws.Range("H113").Value = 3933.3333
$ws.Range("I113").Value = 2533.3333
$ws.Range("J113").Value = 5333.3335
$ws.Range("K113").Value = 2533.3333
$ws.Range("L113").Value = 5333.3335
$ws.Range("M113").Value = -363.3332999999998
$ws.Range("N113").Value = -9673.333500000001

# Row 132
$ws.Range("H132").Value = 1439564.8
$ws.Range("I132").Value = 2978320
$ws.Range("J132").Value = 3393.2
$ws.Range("K132").Value = 8934960
$ws.Range("L132").Value = 10179.6
$ws.Range("M132").Value = -8932430
$ws.Range("N132").Value = -15239.6


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2959.6667
$ws.Range("I7").Value = 2214.5
$ws.Range("J7").Value = 4450
$ws.Range("K7").Value = 2214.5
$ws.Range("L7").Value = 4450
$ws.Range("M7").Value = -2102.5
$ws.Range("N7").Value = -4674

# Row 40
$ws.Range("H40").Value = 5402.636
$ws.Range("I40").Value = 3765.8
$ws.Range("J40").Value = 6766.6665
$ws.Range("K40").Value = 3765.8
$ws.Range("L40").Value = 6766.6665
$ws.Range("M40").Value = -3629.8
$ws.Range("N40").Value = -7038.6665

# Row 126
$ws.Range("H126").Value = 2959.6667
$ws.Range("I126").Value = 2214.5
$ws.Range("J126").Value = 4450
$ws.Range("K126").Value = 6643.5
$ws.Range("L126").Value = 13350
$ws.Range("M126").Value = -4173.5
$ws.Range("N126").Value = -18290

# Row 132
$ws.Range("H132").Value = 66670468
$ws.Range("I132").Value = 83336830
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 250010490
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -250007960
$ws.Range("N132").Value = -20058.0005

# Row 139
$ws.Range("H139").Value = 31982.857
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 31982.857
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 31982.857
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -42262.857


$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 17500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 17500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 17500
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -18748

# Row 66
$ws.Range("H66").Value = 17500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 17500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 52500
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -58740

# Row 82
$ws.Range("H82").Value = 18000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 18000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 18000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -18766

# Row 85
$ws.Range("H85").Value = 18000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 18000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 18000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -20652

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

# Row 98
$ws.Range("H98").Value = 39500
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 39500
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 39500
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -45490

# Row 132
$ws.Range("H132").Value = 3652.2942
$ws.Range("I132").Value = 2924.25
$ws.Range("J132").Value = 5399.6
$ws.Range("K132").Value = 8772.75
$ws.Range("L132").Value = 16198.8
$ws.Range("M132").Value = -6242.75
$ws.Range("N132").Value = -21258.8

# Row 139
$ws.Range("H139").Value = 34112.855
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34112.855
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34112.855
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -44392.855

